# Scheduled-runner style refresh of market-derived columns (H:N) across the
# per-job worksheets. Only the literal values fetched/recomputed by the
# runner change; everything else (labels, dates, item ids, formatting) is
# left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1197458.1
$ws.Range("I132").Value = 2473.9143
$ws.Range("K132").Value = 7421.742899999999
$ws.Range("M132").Value = -4891.742899999999
$ws.Range("H137").Value = 3335244.5
$ws.Range("I137").Value = 5001578
$ws.Range("J137").Value = 2578
$ws.Range("K137").Value = 15004734
$ws.Range("L137").Value = 7734
$ws.Range("M137").Value = -15002184
$ws.Range("N137").Value = -12834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 41750972
$ws.Range("I61").Value = 62563468
$ws.Range("J61").Value = 125978.5
$ws.Range("K61").Value = 62563468
$ws.Range("L61").Value = 125978.5
$ws.Range("M61").Value = -62563256
$ws.Range("N61").Value = -126402.5
$ws.Range("H74").Value = 8130112.5
$ws.Range("I74").Value = 12550857
$ws.Range("J74").Value = 92395.37
$ws.Range("K74").Value = 12550857
$ws.Range("L74").Value = 92395.37
$ws.Range("M74").Value = -12549983
$ws.Range("N74").Value = -94143.37
$ws.Range("H77").Value = 8130112.5
$ws.Range("I77").Value = 12550857
$ws.Range("J77").Value = 92395.37
$ws.Range("K77").Value = 62754285
$ws.Range("L77").Value = 461976.85
$ws.Range("M77").Value = -62749917
$ws.Range("N77").Value = -470712.85
$ws.Range("H132").Value = 40129.152
$ws.Range("I132").Value = 27076.795
$ws.Range("J132").Value = 79286.234
$ws.Range("K132").Value = 81230.38499999999
$ws.Range("L132").Value = 237858.702
$ws.Range("M132").Value = -78700.38499999999
$ws.Range("N132").Value = -242918.702
$ws.Range("H136").Value = 41750972
$ws.Range("I136").Value = 62563468
$ws.Range("J136").Value = 125978.5
$ws.Range("K136").Value = 187690404
$ws.Range("L136").Value = 377935.5
$ws.Range("M136").Value = -187687854
$ws.Range("N136").Value = -383035.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4214.2856
$ws.Range("I134").Value = 3749.3635
$ws.Range("J134").Value = 5919
$ws.Range("K134").Value = 11248.0905
$ws.Range("L134").Value = 17757
$ws.Range("M134").Value = -8713.0905
$ws.Range("N134").Value = -22827

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6792.212
$ws.Range("I31").Value = 47293.668
$ws.Range("J31").Value = 1205.8046
$ws.Range("K31").Value = 47293.668
$ws.Range("L31").Value = 1205.8046
$ws.Range("M31").Value = -46998.668
$ws.Range("N31").Value = -1795.8046
$ws.Range("H34").Value = 6792.212
$ws.Range("I34").Value = 47293.668
$ws.Range("J34").Value = 1205.8046
$ws.Range("K34").Value = 47293.668
$ws.Range("L34").Value = 1205.8046
$ws.Range("M34").Value = -47091.668
$ws.Range("N34").Value = -1609.8046
$ws.Range("H58").Value = 23184410
$ws.Range("I58").Value = 29145724
$ws.Range("J58").Value = 1522
$ws.Range("K58").Value = 29145724
$ws.Range("L58").Value = 1522
$ws.Range("M58").Value = -29145521
$ws.Range("N58").Value = -1928
$ws.Range("H127").Value = 31944.445
$ws.Range("J127").Value = 31944.445
$ws.Range("L127").Value = 31944.445
$ws.Range("N127").Value = -41864.445
$ws.Range("H132").Value = 33322.062
$ws.Range("I132").Value = 1962.25
$ws.Range("J132").Value = 127401.5
$ws.Range("K132").Value = 5886.75
$ws.Range("L132").Value = 382204.5
$ws.Range("M132").Value = -3356.75
$ws.Range("N132").Value = -387264.5
$ws.Range("H134").Value = 33672.47
$ws.Range("I134").Value = 1852.2727
$ws.Range("J134").Value = 92009.5
$ws.Range("K134").Value = 5556.8181
$ws.Range("L134").Value = 276028.5
$ws.Range("M134").Value = -3021.8181
$ws.Range("N134").Value = -281098.5
$ws.Range("H136").Value = 23184410
$ws.Range("I136").Value = 29145724
$ws.Range("J136").Value = 1522
$ws.Range("K136").Value = 87437172
$ws.Range("L136").Value = 4566
$ws.Range("M136").Value = -87434622
$ws.Range("N136").Value = -9666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 18033.701
$ws.Range("I5").Value = 37437.184
$ws.Range("J5").Value = 570.56665
$ws.Range("K5").Value = 112311.552
$ws.Range("L5").Value = 1711.69995
$ws.Range("M5").Value = -112199.552
$ws.Range("N5").Value = -1935.69995
$ws.Range("H68").Value = 1122.5
$ws.Range("I68").Value = 533.17645
$ws.Range("J68").Value = 1957.375
$ws.Range("K68").Value = 1599.52935
$ws.Range("L68").Value = 5872.125
$ws.Range("M68").Value = -788.5293500000002
$ws.Range("N68").Value = -7494.125
$ws.Range("H71").Value = 1122.5
$ws.Range("I71").Value = 533.17645
$ws.Range("J71").Value = 1957.375
$ws.Range("K71").Value = 4798.58805
$ws.Range("L71").Value = 17616.375
$ws.Range("M71").Value = -742.5880500000003
$ws.Range("N71").Value = -25728.375
$ws.Range("H107").Value = 829.25757
$ws.Range("I107").Value = 414.94232
$ws.Range("J107").Value = 2368.1428
$ws.Range("K107").Value = 1244.82696
$ws.Range("L107").Value = 7104.428400000001
$ws.Range("M107").Value = 675.1730400000001
$ws.Range("N107").Value = -10944.4284
$ws.Range("H131").Value = 778.6111
$ws.Range("I131").Value = 466.27274
$ws.Range("J131").Value = 916.04
$ws.Range("K131").Value = 1398.81822
$ws.Range("L131").Value = 2748.12
$ws.Range("M131").Value = 3641.18178
$ws.Range("N131").Value = -12828.12
$ws.Range("H135").Value = 18033.701
$ws.Range("I135").Value = 37437.184
$ws.Range("J135").Value = 570.56665
$ws.Range("K135").Value = 336934.656
$ws.Range("L135").Value = 5135.09985
$ws.Range("M135").Value = -334399.656
$ws.Range("N135").Value = -10205.09985

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 865.6842
$ws.Range("I22").Value = 460.2857
$ws.Range("J22").Value = 1102.1666
$ws.Range("K22").Value = 460.2857
$ws.Range("L22").Value = 1102.1666
$ws.Range("M22").Value = -165.2857
$ws.Range("N22").Value = -1692.1666
$ws.Range("H27").Value = 865.6842
$ws.Range("I27").Value = 460.2857
$ws.Range("J27").Value = 1102.1666
$ws.Range("K27").Value = 460.2857
$ws.Range("L27").Value = 1102.1666
$ws.Range("M27").Value = -353.2857
$ws.Range("N27").Value = -1316.1666
$ws.Range("H132").Value = 38470.57
$ws.Range("I132").Value = 1707.0667
$ws.Range("J132").Value = 80890
$ws.Range("K132").Value = 5121.2001
$ws.Range("L132").Value = 242670
$ws.Range("M132").Value = -2591.2001
$ws.Range("N132").Value = -247730
$ws.Range("H136").Value = 90614.55
$ws.Range("I136").Value = 44086.668
$ws.Range("K136").Value = 132260.004
$ws.Range("M136").Value = -129710.004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1288
$ws.Range("I113").Value = 461.3
$ws.Range("J113").Value = 2665.8333
$ws.Range("K113").Value = 1383.9
$ws.Range("L113").Value = 7997.499899999999
$ws.Range("M113").Value = 786.0999999999999
$ws.Range("N113").Value = -12337.4999
$ws.Range("H132").Value = 59267.855
$ws.Range("I132").Value = 43719.293
$ws.Range("K132").Value = 131157.879
$ws.Range("M132").Value = -128627.879
$ws.Range("H136").Value = 56111.51
$ws.Range("I136").Value = 33352.13
$ws.Range("K136").Value = 100056.39
$ws.Range("M136").Value = -97506.38999999998
